$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
